$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replicate the formatting (style) of the last existing data row's date
# column (A229) onto the new date cells before writing their values, so
# the new rows inherit the same border/font/alignment/number-format as
# the rest of column A.
$ws.Range("A229").Copy($ws.Range("A230:A233"))

$data = @(
    @(230, 44304, 7,  35, 194.6823895872733),
    @(231, 44305, 4,  33, 183.5576816108577),
    @(232, 44306, 11, 40, 222.4941595283124),
    @(233, 44307, 0,  39, 216.9318055401046)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}
